# Updates coinranking.com "cryptos" sheet data to the latest scrape.
# Mirrors the GitHub Actions job that refreshes Coin/Link/Price/Volume(1h)/Hora
# columns on a schedule (commit: "Updated symbol list ... with GitHub Actions").
#
# Values in the Price (D), Volume(1h) (E) and Hora (G) columns are stored as
# literal text in the source sheet (e.g. "0.07700", "1,904.57%", "10"), so we
# prefix them with a leading apostrophe -- exactly like typing into Excel -- to
# stop the auto-detection from converting them into numbers/percentages/dates
# and silently dropping significant trailing zeros, the "%" sign, thousands
# separators, etc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'308.76"
$ws.Range("E2").Value = "'0.03%"
$ws.Range("G2").Value = "'10"

# Row 3
$ws.Range("D3").Value = "'41.18"
$ws.Range("E3").Value = "'1.07%"
$ws.Range("G3").Value = "'10"

# Row 4
$ws.Range("D4").Value = "'5.184"
$ws.Range("E4").Value = "'1.22%"
$ws.Range("G4").Value = "'10"

# Row 5
$ws.Range("D5").Value = "'0.07694"
$ws.Range("E5").Value = "'0.86%"
$ws.Range("G5").Value = "'10"

# Row 6
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.635"
$ws.Range("E6").Value = "'1.64%"
$ws.Range("G6").Value = "'10"

# Row 7
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9147"
$ws.Range("E7").Value = "'1.19%"
$ws.Range("G7").Value = "'10"

# Row 8
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.432"
$ws.Range("E8").Value = "'-1.50%"
$ws.Range("G8").Value = "'10"

# Row 9
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1230"
$ws.Range("E9").Value = "'10.66%"
$ws.Range("G9").Value = "'10"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1828"
$ws.Range("E10").Value = "'2.42%"
$ws.Range("G10").Value = "'10"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09187"
$ws.Range("E11").Value = "'0.58%"
$ws.Range("G11").Value = "'10"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04226"
$ws.Range("E12").Value = "'0.33%"
$ws.Range("G12").Value = "'10"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1050"
$ws.Range("E13").Value = "'-0.11%"
$ws.Range("G13").Value = "'10"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001257"
$ws.Range("E14").Value = "'0.44%"
$ws.Range("G14").Value = "'10"

# Row 15
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005788"
$ws.Range("E15").Value = "'1.47%"
$ws.Range("G15").Value = "'10"

# Row 16
$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").Value = "'0.007509"
$ws.Range("E16").Value = "'1,904.57%"
$ws.Range("G16").Value = "'10"

# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.344"
$ws.Range("E17").Value = "'-0.22%"
$ws.Range("G17").Value = "'10"

# Row 18
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.308"
$ws.Range("E18").Value = "'1.44%"
$ws.Range("G18").Value = "'10"

# Row 19
$ws.Range("E19").Value = "'1.32%"
$ws.Range("G19").Value = "'10"

# Row 20
$ws.Range("D20").Value = "'7.358"
$ws.Range("E20").Value = "'11.23%"
$ws.Range("G20").Value = "'10"

# Row 21
$ws.Range("D21").Value = "'0.1381"
$ws.Range("E21").Value = "'1.59%"
$ws.Range("G21").Value = "'10"

# Row 22
$ws.Range("D22").Value = "'0.2707"
$ws.Range("E22").Value = "'-2.93%"
$ws.Range("G22").Value = "'10"

# Row 23
$ws.Range("D23").Value = "'0.04023"
$ws.Range("E23").Value = "'-1.37%"
$ws.Range("G23").Value = "'10"

# Row 24
$ws.Range("D24").Value = "'0.001262"
$ws.Range("E24").Value = "'2.71%"
$ws.Range("G24").Value = "'10"

# Row 25
$ws.Range("D25").Value = "'0.004267"
$ws.Range("E25").Value = "'6.13%"
$ws.Range("G25").Value = "'10"

# Row 26
$ws.Range("D26").Value = "'0.0001300"
$ws.Range("E26").Value = "'-0.05%"
$ws.Range("G26").Value = "'10"

# Row 27
$ws.Range("G27").Value = "'10"

# Row 28
$ws.Range("G28").Value = "'10"

# Row 29
$ws.Range("G29").Value = "'10"

# Row 30
$ws.Range("G30").Value = "'10"

# Row 31
$ws.Range("G31").Value = "'10"

# Row 32
$ws.Range("G32").Value = "'10"

# Row 33
$ws.Range("G33").Value = "'10"

# Row 34
$ws.Range("G34").Value = "'10"

# Row 35
$ws.Range("G35").Value = "'10"

# Row 36
$ws.Range("G36").Value = "'10"

# Row 37
$ws.Range("G37").Value = "'10"

# Row 38
$ws.Range("D38").Value = "'0.02519"
$ws.Range("E38").Value = "'4.56%"
$ws.Range("G38").Value = "'10"

# Row 39
$ws.Range("D39").Value = "'0.05312"
$ws.Range("E39").Value = "'2.56%"
$ws.Range("G39").Value = "'10"

# Row 40
$ws.Range("D40").Value = "'0.007843"
$ws.Range("E40").Value = "'0.93%"
$ws.Range("G40").Value = "'10"

# Row 41
$ws.Range("D41").Value = "'0.1314"
$ws.Range("E41").Value = "'0.93%"
$ws.Range("G41").Value = "'10"

# Row 42
$ws.Range("D42").Value = "'0.006662"
$ws.Range("E42").Value = "'-5.46%"
$ws.Range("G42").Value = "'10"

# Row 43
$ws.Range("D43").Value = "'0.001860"
$ws.Range("E43").Value = "'-4.66%"
$ws.Range("G43").Value = "'10"

# Row 44
$ws.Range("D44").Value = "'0.007982"
$ws.Range("E44").Value = "'0.63%"
$ws.Range("G44").Value = "'10"

# Row 45
$ws.Range("D45").Value = "'0.3069"
$ws.Range("E45").Value = "'-0.47%"
$ws.Range("G45").Value = "'10"

# Row 46
$ws.Range("D46").Value = "'0.00006700"
$ws.Range("E46").Value = "'-4.10%"
$ws.Range("G46").Value = "'10"

# Row 47
$ws.Range("G47").Value = "'10"

# Row 48
$ws.Range("D48").Value = "'0.3012"
$ws.Range("E48").Value = "'854.23%"
$ws.Range("G48").Value = "'10"

# Row 49
$ws.Range("G49").Value = "'10"

# Row 50
$ws.Range("G50").Value = "'10"

# Row 51
$ws.Range("G51").Value = "'10"
